$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tolerances")

# Set the value of C9 (Maximum accepted deviation from COG (metres)) to 1
$ws.Range("C9").Value = 1

# Update the active cell selection to D19 (as recorded in the workbook view state)
$ws.Range("D19").Select()
